$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 used to be "=W22" (a formula pointing at a country further down the
# lookup table). Replace it with a plain value pointing at the first
# lookup-table row (W2, "Altafia") so the VLOOKUPs below resolve to that
# country's numbers (a cheaper / "upgradable" save tier).
$ws.Range("A1").Formula = ""
$ws.Range("A1").Value = "Altafia"

# Record the chosen income value in a new helper cell, matching the
# formatting used elsewhere in the sheet (#,##0).
$ws.Range("X1").Value = 749923
$ws.Range("X1").NumberFormat = "#,##0"

# Unhide the helper columns J:S (10-19) that hold the tier math; column T
# stays hidden.
$ws.Range("J1:S1").EntireColumn.Hidden = $false

# Update the view so it matches where the user left off scrolling/selecting.
$ws.Range("F5").Select()
$ws.Application.ActiveWindow.ScrollColumn = 4
